# Esboço da função fight()
# Fill in the "Resultado" (H) / "Comentário" (I) status columns for the
# test-plan rows that cover the new fight() function draft.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: mark as "esperado" / "-" -------------------------------------
# H4 changes from "-" to "esperado" and needs the green "esperado" look
# (same formatting already used in H3), so copy H3's format over first.
$ws.Range("H3").Copy()
$ws.Range("H4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H4").Value = "esperado"
$ws.Range("I4").Value = "-"

# --- Row 5: fill in previously empty status cells -------------------------
$ws.Range("H3").Copy()
$ws.Range("H5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H5").Value = "esperado"
$ws.Range("I5").Value = "-"

# --- Row 6: fill in previously empty status cells (format already right) --
$ws.Range("H6").Value = "-"
$ws.Range("I6").Value = "-"

$excel.CutCopyMode = 0

# --- Column widths nudge slightly narrower after the edits ---------------
$ws.Columns.Item(1).ColumnWidth = 9.333333333333334
$ws.Columns.Item(2).ColumnWidth = 13
$ws.Columns.Item(3).ColumnWidth = 13.666666666666666
$ws.Columns.Item(4).ColumnWidth = 14.833333333333334
$ws.Columns.Item(5).ColumnWidth = 15.666666666666666
$ws.Columns.Item(6).ColumnWidth = 34.5
$ws.Columns.Item(7).ColumnWidth = 30.5
$ws.Columns.Item(8).ColumnWidth = 11
$ws.Columns.Item(9).ColumnWidth = 17.5

# --- Update the active view/selection -------------------------------------
$ws.Range("A3").Select()
